# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the 385087e5-... file's handback row.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-24 02:47:53"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-24 02:47:48"
$zhcn.Range("K3").Value = "2016-08-24 02:48:12"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-24 02:47:53"
$dede.Range("K3").Value = "2016-08-24 02:48:19"
